$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update dSF column (F) values for specific rows
$ws.Range("F2").Value = -4
$ws.Range("F3").Value = -2
$ws.Range("F6").Value = 0
$ws.Range("F8").Value = -4
